$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("06_API一覧")

# Fix API ID numbering: the "publish" row (API-031 data duplicated at the
# bottom of the app-scope block) is moved back to directly follow API-030
# (row 32), and every row that was API-032..API-065 shifts down by one row.
# Column A (No) already holds the correct sequential 1..86 values and the
# admin-scope rows (67+) are untouched, so only columns B..M for rows 32..66
# need to be rewritten.

# Row 32
$ws.Range("B32").Value = 'API-031'
$ws.Range("C32").Value = '/api/v1/instructor/courses/{courseId}/publish'
$ws.Range("D32").Value = 'POST'
$ws.Range("E32").Value = 'コース公開（講師）'
$ws.Range("F32").Value = 'app'
$ws.Range("G32").Value = 'instructor_owner'
$ws.Range("H32").Value = 'AUDIT_LOG'
$ws.Range("I32").Value = '-'
$ws.Range("J32").Value = 'CourseDetailView'
$ws.Range("K32").Value = '200, 403, 423'
$ws.Range("L32").Value = 'v1.2'
$ws.Range("M32").Value = 'ownerUserId一致の講師がコースを公開。statusをactiveに変更。operator作成コースの場合、CourseMember.role=instructor→instructor_ownerへ昇格（委譲完了）。'

# Row 33
$ws.Range("B33").Value = 'API-032'
$ws.Range("C33").Value = '/api/v1/instructor/analytics/overview'
$ws.Range("D33").Value = 'GET'
$ws.Range("E33").Value = '売上分析サマリ'
$ws.Range("F33").Value = 'app'
$ws.Range("G33").Value = 'instructor'
$ws.Range("H33").Value = '-'
$ws.Range("I33").Value = '-'
$ws.Range("J33").Value = 'GenericListResponse'
$ws.Range("K33").Value = '200'
$ws.Range("L33").Value = 'KEEP'
$ws.Range("M33").Value = '期間別の収益、成約率、受講生離脱ポイント等の分析サマリーを取得。'

# Row 34
$ws.Range("B34").Value = 'API-033'
$ws.Range("C34").Value = '/api/v1/instructor/settings/payout'
$ws.Range("D34").Value = 'PUT'
$ws.Range("E34").Value = '振込先銀行設定'
$ws.Range("F34").Value = 'app'
$ws.Range("G34").Value = 'instructor_owner'
$ws.Range("H34").Value = '-'
$ws.Range("I34").Value = 'GenericWriteRequest'
$ws.Range("J34").Value = 'SuccessResponse'
$ws.Range("K34").Value = '200'
$ws.Range("L34").Value = 'KEEP'
$ws.Range("M34").Value = '講師への報酬振込先口座情報をStripe Connect等と連携して設定・更新。'

# Row 35
$ws.Range("B35").Value = 'API-034'
$ws.Range("C35").Value = '/api/v1/instructor/courses/{courseId}/syllabus'
$ws.Range("D35").Value = 'GET'
$ws.Range("E35").Value = 'シラバス構造取得'
$ws.Range("F35").Value = 'app'
$ws.Range("G35").Value = 'instructor'
$ws.Range("H35").Value = '-'
$ws.Range("I35").Value = '-'
$ws.Range("J35").Value = 'CourseDetailView'
$ws.Range("K35").Value = '200'
$ws.Range("L35").Value = 'KEEP'
$ws.Range("M35").Value = '編集中の章立てとレッスン構成をツリー形式で取得（image_32dd18の構成用）。'

# Row 36
$ws.Range("B36").Value = 'API-035'
$ws.Range("C36").Value = '/api/v1/instructor/courses/{courseId}/sections'
$ws.Range("D36").Value = 'POST'
$ws.Range("E36").Value = 'セクション追加'
$ws.Range("F36").Value = 'app'
$ws.Range("G36").Value = 'instructor'
$ws.Range("H36").Value = '423_ON_FROZEN'
$ws.Range("I36").Value = 'GenericWriteRequest'
$ws.Range("J36").Value = 'CourseDetailView'
$ws.Range("K36").Value = '201, 423'
$ws.Range("L36").Value = 'KEEP'
$ws.Range("M36").Value = 'カリキュラム内に新しい「章」を追加。423チェック対象。'

# Row 37
$ws.Range("B37").Value = 'API-036'
$ws.Range("C37").Value = '/api/v1/instructor/sections/{sectionId}'
$ws.Range("D37").Value = 'PUT'
$ws.Range("E37").Value = 'セクション編集'
$ws.Range("F37").Value = 'app'
$ws.Range("G37").Value = 'instructor'
$ws.Range("H37").Value = '423_ON_FROZEN'
$ws.Range("I37").Value = 'GenericWriteRequest'
$ws.Range("J37").Value = 'SuccessResponse'
$ws.Range("K37").Value = '200, 423'
$ws.Range("L37").Value = 'KEEP'
$ws.Range("M37").Value = 'セクション名の変更、表示順の並び替え、および削除。'

# Row 38
$ws.Range("B38").Value = 'API-037'
$ws.Range("C38").Value = '/api/v1/instructor/sections/{sectionId}'
$ws.Range("D38").Value = 'DELETE'
$ws.Range("E38").Value = 'セクション削除'
$ws.Range("F38").Value = 'app'
$ws.Range("G38").Value = 'instructor'
$ws.Range("H38").Value = '423_ON_FROZEN'
$ws.Range("I38").Value = '-'
$ws.Range("J38").Value = 'SuccessResponse'
$ws.Range("K38").Value = '200, 423'
$ws.Range("L38").Value = 'KEEP'
$ws.Range("M38").Value = '指定されたセクションと、配下のレッスン紐付けを削除（論理/物理）。'

# Row 39
$ws.Range("B39").Value = 'API-038'
$ws.Range("C39").Value = '/api/v1/instructor/sections/{sectionId}/lessons'
$ws.Range("D39").Value = 'POST'
$ws.Range("E39").Value = 'レッスン作成'
$ws.Range("F39").Value = 'app'
$ws.Range("G39").Value = 'instructor'
$ws.Range("H39").Value = '423_ON_FROZEN'
$ws.Range("I39").Value = 'GenericWriteRequest'
$ws.Range("J39").Value = 'SuccessResponse'
$ws.Range("K39").Value = '201, 423'
$ws.Range("L39").Value = 'KEEP'
$ws.Range("M39").Value = '指定セクションに新規レッスンを追加（座学/動画/課題/ライブ）（image_32dcf8反映）。'

# Row 40
$ws.Range("B40").Value = 'API-039'
$ws.Range("C40").Value = '/api/v1/instructor/lessons/{lessonId}'
$ws.Range("D40").Value = 'GET'
$ws.Range("E40").Value = 'レッスン詳細取得'
$ws.Range("F40").Value = 'app'
$ws.Range("G40").Value = 'instructor'
$ws.Range("H40").Value = '-'
$ws.Range("I40").Value = '-'
$ws.Range("J40").Value = 'GenericDetailView'
$ws.Range("K40").Value = '200'
$ws.Range("L40").Value = 'KEEP'
$ws.Range("M40").Value = '編集画面用。動画URL、リッチテキスト本文、課題要件、Drip設定の詳細。'

# Row 41
$ws.Range("B41").Value = 'API-040'
$ws.Range("C41").Value = '/api/v1/instructor/lessons/{lessonId}'
$ws.Range("D41").Value = 'PUT'
$ws.Range("E41").Value = 'レッスン編集'
$ws.Range("F41").Value = 'app'
$ws.Range("G41").Value = 'instructor'
$ws.Range("H41").Value = '423_ON_FROZEN'
$ws.Range("I41").Value = 'GenericWriteRequest'
$ws.Range("J41").Value = 'SuccessResponse'
$ws.Range("K41").Value = '200, 423'
$ws.Range("L41").Value = 'KEEP'
$ws.Range("M41").Value = 'コンテンツの修正、Drip(解禁日)、先行ロック条件の保存。423チェック対象。'

# Row 42
$ws.Range("B42").Value = 'API-041'
$ws.Range("C42").Value = '/api/v1/instructor/lessons/{lessonId}'
$ws.Range("D42").Value = 'DELETE'
$ws.Range("E42").Value = 'レッスン削除'
$ws.Range("F42").Value = 'app'
$ws.Range("G42").Value = 'instructor'
$ws.Range("H42").Value = '423_ON_FROZEN'
$ws.Range("I42").Value = '-'
$ws.Range("J42").Value = 'SuccessResponse'
$ws.Range("K42").Value = '200, 423'
$ws.Range("L42").Value = 'KEEP'
$ws.Range("M42").Value = 'レッスンをカリキュラムから削除。'

# Row 43
$ws.Range("B43").Value = 'API-042'
$ws.Range("C43").Value = '/api/v1/instructor/lessons/{lessonId}/live'
$ws.Range("D43").Value = 'PUT'
$ws.Range("E43").Value = 'ライブ配信設定'
$ws.Range("F43").Value = 'app'
$ws.Range("G43").Value = 'instructor'
$ws.Range("H43").Value = '423_ON_FROZEN'
$ws.Range("I43").Value = 'GenericWriteRequest'
$ws.Range("J43").Value = 'SuccessResponse'
$ws.Range("K43").Value = '200, 423'
$ws.Range("L43").Value = 'KEEP'
$ws.Range("M43").Value = 'ZoomやYouTube LiveのURL、および配信日時の設定（カレンダー連動）。'

# Row 44
$ws.Range("B44").Value = 'API-043'
$ws.Range("C44").Value = '/api/v1/instructor/courses/{courseId}/submissions'
$ws.Range("D44").Value = 'GET'
$ws.Range("E44").Value = '提出一覧取得'
$ws.Range("F44").Value = 'app'
$ws.Range("G44").Value = 'instructor'
$ws.Range("H44").Value = '-'
$ws.Range("I44").Value = '-'
$ws.Range("J44").Value = 'CourseDetailView'
$ws.Range("K44").Value = '200'
$ws.Range("L44").Value = 'KEEP'
$ws.Range("M44").Value = '受講者別の課題進捗、未確認提出物、最終提出日時の一覧（image_32dcd7反映）。'

# Row 45
$ws.Range("B45").Value = 'API-044'
$ws.Range("C45").Value = '/api/v1/instructor/submissions/{submissionId}/evaluation'
$ws.Range("D45").Value = 'PATCH'
$ws.Range("E45").Value = '評価・採点実行'
$ws.Range("F45").Value = 'app'
$ws.Range("G45").Value = 'instructor, assistant'
$ws.Range("H45").Value = 'THREAD_REPLY(AUTO)'
$ws.Range("I45").Value = 'GenericWriteRequest'
$ws.Range("J45").Value = 'SubmissionView'
$ws.Range("K45").Value = '200'
$ws.Range("L45").Value = 'KEEP'
$ws.Range("M45").Value = '提出物への合否、講評を送信。合格時は自動で次の教材を解禁。フィードバックはスレッドへ。'

# Row 46
$ws.Range("B46").Value = 'API-045'
$ws.Range("C46").Value = '/api/v1/instructor/courses/{courseId}/members'
$ws.Range("D46").Value = 'GET'
$ws.Range("E46").Value = '受講者名簿取得'
$ws.Range("F46").Value = 'app'
$ws.Range("G46").Value = 'instructor'
$ws.Range("H46").Value = '-'
$ws.Range("I46").Value = '-'
$ws.Range("J46").Value = 'CourseDetailView'
$ws.Range("K46").Value = '200'
$ws.Range("L46").Value = 'KEEP'
$ws.Range("M46").Value = '全参加ユーザーの属性、現在の進捗状況、メールアドレス等の名簿取得。'

# Row 47
$ws.Range("B47").Value = 'API-046'
$ws.Range("C47").Value = '/api/v1/instructor/courses/{courseId}/members/{userId}/role'
$ws.Range("D47").Value = 'PATCH'
$ws.Range("E47").Value = '講座内ロール変更'
$ws.Range("F47").Value = 'app'
$ws.Range("G47").Value = 'instructor_owner'
$ws.Range("H47").Value = '423_ON_FROZEN'
$ws.Range("I47").Value = 'GenericWriteRequest'
$ws.Range("J47").Value = 'CourseDetailView'
$ws.Range("K47").Value = '200, 423'
$ws.Range("L47").Value = 'KEEP'
$ws.Range("M47").Value = '特定ユーザーの役割（講師、アシスタント、受講生）を動的に切り替え。'

# Row 48
$ws.Range("B48").Value = 'API-047'
$ws.Range("C48").Value = '/api/v1/instructor/courses/{courseId}/members/{userId}/revoke'
$ws.Range("D48").Value = 'POST'
$ws.Range("E48").Value = '受講権限剥奪'
$ws.Range("F48").Value = 'app'
$ws.Range("G48").Value = 'instructor_owner'
$ws.Range("H48").Value = '423_ON_FROZEN'
$ws.Range("I48").Value = 'GenericWriteRequest'
$ws.Range("J48").Value = 'CourseDetailView'
$ws.Range("K48").Value = '201, 423'
$ws.Range("L48").Value = 'KEEP'
$ws.Range("M48").Value = '特定ユーザーの受講権限を剥奪（revoked）し、コンテンツへのアクセスを遮断。'

# Row 49
$ws.Range("B49").Value = 'API-048'
$ws.Range("C49").Value = '/api/v1/instructor/courses/{courseId}/export'
$ws.Range("D49").Value = 'GET'
$ws.Range("E49").Value = '受講者CSV出力'
$ws.Range("F49").Value = 'app'
$ws.Range("G49").Value = 'instructor_owner'
$ws.Range("H49").Value = '-'
$ws.Range("I49").Value = '-'
$ws.Range("J49").Value = 'CourseDetailView'
$ws.Range("K49").Value = '200'
$ws.Range("L49").Value = 'KEEP'
$ws.Range("M49").Value = '受講者名簿と学習進捗データをCSV形式でエクスポートする。'

# Row 50
$ws.Range("B50").Value = 'API-049'
$ws.Range("C50").Value = '/api/v1/courses/{courseId}/channels'
$ws.Range("D50").Value = 'GET'
$ws.Range("E50").Value = 'チャンネル一覧取得'
$ws.Range("F50").Value = 'app'
$ws.Range("G50").Value = 'all_in_course'
$ws.Range("H50").Value = '-'
$ws.Range("I50").Value = '-'
$ws.Range("J50").Value = 'CourseChannelListResponse'
$ws.Range("K50").Value = '200'
$ws.Range("L50").Value = 'KEEP'
$ws.Range("M50").Value = 'サイドメニュー構築用。general, announcement等のチャンネル種別を含む。'

# Row 51
$ws.Range("B51").Value = 'API-050'
$ws.Range("C51").Value = '/api/v1/courses/{courseId}/channels'
$ws.Range("D51").Value = 'POST'
$ws.Range("E51").Value = 'チャンネル作成'
$ws.Range("F51").Value = 'app'
$ws.Range("G51").Value = 'instructor_owner'
$ws.Range("H51").Value = '423_ON_FROZEN'
$ws.Range("I51").Value = 'CourseChannelCreateRequest'
$ws.Range("J51").Value = 'CourseDetailView'
$ws.Range("K51").Value = '201, 423'
$ws.Range("L51").Value = 'KEEP'
$ws.Range("M51").Value = 'コース内に新しいカスタムチャンネルを追加。'

# Row 52
$ws.Range("B52").Value = 'API-051'
$ws.Range("C52").Value = '/api/v1/channels/{channelId}'
$ws.Range("D52").Value = 'PUT'
$ws.Range("E52").Value = 'チャンネル編集'
$ws.Range("F52").Value = 'app'
$ws.Range("G52").Value = 'instructor_owner'
$ws.Range("H52").Value = '423_ON_FROZEN'
$ws.Range("I52").Value = 'GenericWriteRequest'
$ws.Range("J52").Value = 'CourseChannelDetailView'
$ws.Range("K52").Value = '200, 423'
$ws.Range("L52").Value = 'KEEP'
$ws.Range("M52").Value = 'チャンネル名、説明文、閲覧制限、アーカイブ状態の更新。'

# Row 53
$ws.Range("B53").Value = 'API-052'
$ws.Range("C53").Value = '/api/v1/channels/{channelId}'
$ws.Range("D53").Value = 'DELETE'
$ws.Range("E53").Value = 'チャンネル削除'
$ws.Range("F53").Value = 'app'
$ws.Range("G53").Value = 'instructor_owner'
$ws.Range("H53").Value = '423_ON_FROZEN'
$ws.Range("I53").Value = '-'
$ws.Range("J53").Value = 'SuccessResponse'
$ws.Range("K53").Value = '200, 423'
$ws.Range("L53").Value = 'KEEP'
$ws.Range("M53").Value = 'チャンネルを論理削除。過去ログは監査用に保持。'

# Row 54
$ws.Range("B54").Value = 'API-053'
$ws.Range("C54").Value = '/api/v1/channels/{channelId}/messages'
$ws.Range("D54").Value = 'GET'
$ws.Range("E54").Value = 'メッセージ履歴取得'
$ws.Range("F54").Value = 'app'
$ws.Range("G54").Value = 'all_in_course'
$ws.Range("H54").Value = '閲覧のみ可'
$ws.Range("I54").Value = '-'
$ws.Range("J54").Value = 'MessageListResponse'
$ws.Range("K54").Value = '200'
$ws.Range("L54").Value = 'KEEP'
$ws.Range("M54").Value = 'スレッド親（ルート）メッセージの一覧取得。課題の相互閲覧可（image_32d9b1反映）。'

# Row 55
$ws.Range("B55").Value = 'API-054'
$ws.Range("C55").Value = '/api/v1/channels/{channelId}/messages'
$ws.Range("D55").Value = 'POST'
$ws.Range("E55").Value = 'メッセージ投稿'
$ws.Range("F55").Value = 'app'
$ws.Range("G55").Value = 'all_in_course'
$ws.Range("H55").Value = 'threads_only(AUTO)'
$ws.Range("I55").Value = 'CourseMessageCreateRequest'
$ws.Range("J55").Value = 'CourseChannelDetailView'
$ws.Range("K55").Value = '201'
$ws.Range("L55").Value = 'KEEP'
$ws.Range("M55").Value = 'チャンネルへの新規投稿。親メッセージとして保存。threads_onlyを強制適用。'

# Row 56
$ws.Range("B56").Value = 'API-055'
$ws.Range("C56").Value = '/api/v1/messages/{messageId}/thread'
$ws.Range("D56").Value = 'GET'
$ws.Range("E56").Value = 'スレッド詳細取得'
$ws.Range("F56").Value = 'app'
$ws.Range("G56").Value = 'all_in_course'
$ws.Range("H56").Value = '-'
$ws.Range("I56").Value = '-'
$ws.Range("J56").Value = 'GenericListResponse'
$ws.Range("K56").Value = '200'
$ws.Range("L56").Value = 'KEEP'
$ws.Range("M56").Value = '特定メッセージに紐づく返信メッセージ一覧（スレッドビュー）を全件取得。'

# Row 57
$ws.Range("B57").Value = 'API-056'
$ws.Range("C57").Value = '/api/v1/messages/{messageId}/replies'
$ws.Range("D57").Value = 'POST'
$ws.Range("E57").Value = 'スレッド返信投稿'
$ws.Range("F57").Value = 'app'
$ws.Range("G57").Value = 'all_in_course'
$ws.Range("H57").Value = 'THREAD_REPLY'
$ws.Range("I57").Value = 'GenericWriteRequest'
$ws.Range("J57").Value = 'CourseMessageView'
$ws.Range("K57").Value = '201'
$ws.Range("L57").Value = 'KEEP'
$ws.Range("M57").Value = '返信投稿(image_32d9b1)。announcementタイプは講師以外返信不可。'

# Row 58
$ws.Range("B58").Value = 'API-057'
$ws.Range("C58").Value = '/api/v1/messages/{messageId}'
$ws.Range("D58").Value = 'PATCH'
$ws.Range("E58").Value = 'メッセージ編集'
$ws.Range("F58").Value = 'app'
$ws.Range("G58").Value = 'owner_only'
$ws.Range("H58").Value = '423_ON_FROZEN'
$ws.Range("I58").Value = 'GenericWriteRequest'
$ws.Range("J58").Value = 'CourseMessageView'
$ws.Range("K58").Value = '200, 423'
$ws.Range("L58").Value = 'KEEP'
$ws.Range("M58").Value = '自分の投稿内容を修正。編集履歴を保持。423チェック対象。'

# Row 59
$ws.Range("B59").Value = 'API-058'
$ws.Range("C59").Value = '/api/v1/messages/{messageId}'
$ws.Range("D59").Value = 'DELETE'
$ws.Range("E59").Value = 'メッセージ削除'
$ws.Range("F59").Value = 'app'
$ws.Range("G59").Value = 'owner_only, instructor'
$ws.Range("H59").Value = '423_ON_FROZEN'
$ws.Range("I59").Value = '-'
$ws.Range("J59").Value = 'SuccessResponse'
$ws.Range("K59").Value = '200, 423'
$ws.Range("L59").Value = 'KEEP'
$ws.Range("M59").Value = '投稿の論理削除。返信がある場合は「削除されました」と表示。'

# Row 60
$ws.Range("B60").Value = 'API-059'
$ws.Range("C60").Value = '/api/v1/messages/{messageId}/reactions'
$ws.Range("D60").Value = 'POST'
$ws.Range("E60").Value = 'リアクション追加'
$ws.Range("F60").Value = 'app'
$ws.Range("G60").Value = 'all_in_course'
$ws.Range("H60").Value = '423_ON_FROZEN'
$ws.Range("I60").Value = 'GenericWriteRequest'
$ws.Range("J60").Value = 'CourseMessageView'
$ws.Range("K60").Value = '201, 423'
$ws.Range("L60").Value = 'KEEP'
$ws.Range("M60").Value = 'メッセージに対する絵文字リアクションの付与。423チェック対象。'

# Row 61
$ws.Range("B61").Value = 'API-060'
$ws.Range("C61").Value = '/api/v1/courses/{courseId}/channels/{channelId}/threads'
$ws.Range("D61").Value = 'GET'
$ws.Range("E61").Value = 'スレッド一覧取得（ルートメッセージのみ）'
$ws.Range("F61").Value = 'app'
$ws.Range("G61").Value = 'all_in_course'
$ws.Range("H61").Value = '-'
$ws.Range("I61").Value = '-'
$ws.Range("J61").Value = 'ThreadListResponse'
$ws.Range("K61").Value = '200'
$ws.Range("L61").Value = 'KEEP'
$ws.Range("M61").Value = 'threads_only対応'

# Row 62
$ws.Range("B62").Value = 'API-061'
$ws.Range("C62").Value = '/api/v1/courses/{courseId}/channels/{channelId}/threads'
$ws.Range("D62").Value = 'POST'
$ws.Range("E62").Value = 'スレッド作成（ルート投稿）'
$ws.Range("F62").Value = 'app'
$ws.Range("G62").Value = 'all_in_course'
$ws.Range("H62").Value = '-'
$ws.Range("I62").Value = 'GenericWriteRequest'
$ws.Range("J62").Value = 'CourseDetailView'
$ws.Range("K62").Value = '201'
$ws.Range("L62").Value = 'KEEP'
$ws.Range("M62").Value = 'threadId=NULLで作成'

# Row 63
$ws.Range("B63").Value = 'API-062'
$ws.Range("C63").Value = '/api/v1/courses/{courseId}/channels/{channelId}/threads/{threadId}/messages'
$ws.Range("D63").Value = 'GET'
$ws.Range("E63").Value = 'スレッド内メッセージ一覧（ルート+返信）'
$ws.Range("F63").Value = 'app'
$ws.Range("G63").Value = 'all_in_course'
$ws.Range("H63").Value = '-'
$ws.Range("I63").Value = '-'
$ws.Range("J63").Value = 'MessageListResponse'
$ws.Range("K63").Value = '200'
$ws.Range("L63").Value = 'KEEP'
$ws.Range("M63").Value = 'threadId指定で返信も取得'

# Row 64
$ws.Range("B64").Value = 'API-063'
$ws.Range("C64").Value = '/api/v1/courses/{courseId}/channels/{channelId}/threads/{threadId}/messages'
$ws.Range("D64").Value = 'POST'
$ws.Range("E64").Value = 'スレッド返信'
$ws.Range("F64").Value = 'app'
$ws.Range("G64").Value = 'all_in_course'
$ws.Range("H64").Value = '-'
$ws.Range("I64").Value = 'CourseMessageCreateRequest'
$ws.Range("J64").Value = 'CourseDetailView'
$ws.Range("K64").Value = '201'
$ws.Range("L64").Value = 'KEEP'
$ws.Range("M64").Value = 'threadId必須'

# Row 65
$ws.Range("B65").Value = 'API-064'
$ws.Range("C65").Value = '/api/v1/payments/webhook'
$ws.Range("D65").Value = 'POST'
$ws.Range("E65").Value = 'Stripe Webhook'
$ws.Range("F65").Value = 'app'
$ws.Range("G65").Value = 'public(Stripe)'
$ws.Range("H65").Value = 'Webhook専用'
$ws.Range("I65").Value = 'GenericWriteRequest'
$ws.Range("J65").Value = 'SuccessResponse'
$ws.Range("K65").Value = '200'
$ws.Range("L65").Value = 'KEEP'
$ws.Range("M65").Value = '外部決済完了通知を受け取り、Enrollmentを自動的にactive化する。'

# Row 66
$ws.Range("B66").Value = 'API-065'
$ws.Range("C66").Value = '/api/v1/health'
$ws.Range("D66").Value = 'GET'
$ws.Range("E66").Value = 'ヘルスチェック'
$ws.Range("F66").Value = 'app'
$ws.Range("G66").Value = 'public'
$ws.Range("H66").Value = '-'
$ws.Range("I66").Value = '-'
$ws.Range("J66").Value = 'GenericListResponse'
$ws.Range("K66").Value = '200'
$ws.Range("L66").Value = 'KEEP'
$ws.Range("M66").Value = 'サーバーの死活監視用。DB、Redis、外部サービスの接続確認。'
